$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H62").Value = 2336.3333
$ws.Range("I62").Value = 1729.2727
$ws.Range("J62").Value = 2687.7896
$ws.Range("K62").Value = 1729.2727
$ws.Range("L62").Value = 2687.7896
$ws.Range("M62").Value = -1105.2727
$ws.Range("N62").Value = -3935.7896
$ws.Range("H64").Value = 3837.6924
$ws.Range("I64").Value = 3318
$ws.Range("J64").Value = 4162.5
$ws.Range("K64").Value = 3318
$ws.Range("L64").Value = 4162.5
$ws.Range("M64").Value = -3070
$ws.Range("N64").Value = -4658.5
$ws.Range("H65").Value = 2336.3333
$ws.Range("I65").Value = 1729.2727
$ws.Range("J65").Value = 2687.7896
$ws.Range("K65").Value = 8646.363499999999
$ws.Range("L65").Value = 13438.948
$ws.Range("M65").Value = -5526.363499999999
$ws.Range("N65").Value = -19678.948
$ws.Range("H67").Value = 3837.6924
$ws.Range("I67").Value = 3318
$ws.Range("J67").Value = 4162.5
$ws.Range("K67").Value = 3318
$ws.Range("L67").Value = 4162.5
$ws.Range("M67").Value = -2460
$ws.Range("N67").Value = -5878.5
$ws.Range("H106").Value = 4565.8335
$ws.Range("I106").Value = 5041.4287
$ws.Range("J106").Value = 3900
$ws.Range("K106").Value = 5041.4287
$ws.Range("L106").Value = 3900
$ws.Range("M106").Value = -4410.4287
$ws.Range("N106").Value = -5162
$ws.Range("H138").Value = 1668524.9
$ws.Range("I138").Value = 1166.5518
$ws.Range("J138").Value = 3971067.2
$ws.Range("K138").Value = 3499.6554
$ws.Range("L138").Value = 11913201.6
$ws.Range("M138").Value = 1640.3446
$ws.Range("N138").Value = -11923481.6
$ws.Range("H141").Value = 1701.1228
$ws.Range("I141").Value = 1131.9038
$ws.Range("K141").Value = 3395.7114
$ws.Range("M141").Value = 1784.2886
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3154.2942
$ws.Range("I2").Value = 2387.2856
$ws.Range("J2").Value = 3691.2
$ws.Range("K2").Value = 2387.2856
$ws.Range("L2").Value = 3691.2
$ws.Range("M2").Value = -2274.2856
$ws.Range("N2").Value = -3917.2
$ws.Range("H32").Value = 1589.43
$ws.Range("I32").Value = 1000.8461
$ws.Range("J32").Value = 3676.2273
$ws.Range("K32").Value = 1000.8461
$ws.Range("L32").Value = 3676.2273
$ws.Range("M32").Value = -713.8461
$ws.Range("N32").Value = -4250.2273
$ws.Range("H61").Value = 18906710
$ws.Range("I61").Value = 22245454
$ws.Range("J61").Value = 126280.25
$ws.Range("K61").Value = 22245454
$ws.Range("L61").Value = 126280.25
$ws.Range("M61").Value = -22245242
$ws.Range("N61").Value = -126704.25
$ws.Range("H116").Value = 3154.2942
$ws.Range("I116").Value = 2387.2856
$ws.Range("J116").Value = 3691.2
$ws.Range("K116").Value = 2387.2856
$ws.Range("L116").Value = 3691.2
$ws.Range("M116").Value = -93.28560000000016
$ws.Range("N116").Value = -8279.200000000001
$ws.Range("H132").Value = 60580.6
$ws.Range("I132").Value = 37864.184
$ws.Range("J132").Value = 137248.5
$ws.Range("K132").Value = 113592.552
$ws.Range("L132").Value = 411745.5
$ws.Range("M132").Value = -111062.552
$ws.Range("N132").Value = -416805.5
$ws.Range("H136").Value = 18906710
$ws.Range("I136").Value = 22245454
$ws.Range("J136").Value = 126280.25
$ws.Range("K136").Value = 66736362
$ws.Range("L136").Value = 378840.75
$ws.Range("M136").Value = -66733812
$ws.Range("N136").Value = -383940.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3154.2942
$ws.Range("I3").Value = 2387.2856
$ws.Range("J3").Value = 3691.2
$ws.Range("K3").Value = 2387.2856
$ws.Range("L3").Value = 3691.2
$ws.Range("M3").Value = -2273.2856
$ws.Range("N3").Value = -3919.2
$ws.Range("H80").Value = 276.0909
$ws.Range("I80").Value = 264.375
$ws.Range("J80").Value = 282.7857
$ws.Range("K80").Value = 264.375
$ws.Range("L80").Value = 282.7857
$ws.Range("M80").Value = 733.625
$ws.Range("N80").Value = -2278.7857
$ws.Range("H83").Value = 276.0909
$ws.Range("I83").Value = 264.375
$ws.Range("J83").Value = 282.7857
$ws.Range("K83").Value = 1321.875
$ws.Range("L83").Value = 1413.9285
$ws.Range("M83").Value = 3670.125
$ws.Range("N83").Value = -11397.9285
$ws.Range("H134").Value = 1346
$ws.Range("I134").Value = 856.83673
$ws.Range("J134").Value = 2943.9333
$ws.Range("K134").Value = 2570.51019
$ws.Range("L134").Value = 8831.7999
$ws.Range("M134").Value = -35.51018999999997
$ws.Range("N134").Value = -13901.7999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 90909530
$ws.Range("I22").Value = 166666880
$ws.Range("K22").Value = 166666880
$ws.Range("M22").Value = -166666530
$ws.Range("H31").Value = 3365.2354
$ws.Range("I31").Value = 1366.9565
$ws.Range("J31").Value = 7543.4546
$ws.Range("K31").Value = 1366.9565
$ws.Range("L31").Value = 7543.4546
$ws.Range("M31").Value = -1071.9565
$ws.Range("N31").Value = -8133.4546
$ws.Range("H34").Value = 3365.2354
$ws.Range("I34").Value = 1366.9565
$ws.Range("J34").Value = 7543.4546
$ws.Range("K34").Value = 1366.9565
$ws.Range("L34").Value = 7543.4546
$ws.Range("M34").Value = -1164.9565
$ws.Range("N34").Value = -7947.4546
$ws.Range("H62").Value = 3135.3333
$ws.Range("I62").Value = 2400
$ws.Range("J62").Value = 3503
$ws.Range("K62").Value = 2400
$ws.Range("L62").Value = 3503
$ws.Range("M62").Value = -1776
$ws.Range("N62").Value = -4751
$ws.Range("H65").Value = 3135.3333
$ws.Range("I65").Value = 2400
$ws.Range("J65").Value = 3503
$ws.Range("K65").Value = 12000
$ws.Range("L65").Value = 17515
$ws.Range("M65").Value = -8880
$ws.Range("N65").Value = -23755
$ws.Range("H86").Value = 2447.9524
$ws.Range("I86").Value = 2430.7144
$ws.Range("J86").Value = 2482.4285
$ws.Range("K86").Value = 2430.7144
$ws.Range("L86").Value = 2482.4285
$ws.Range("M86").Value = -1307.7144
$ws.Range("N86").Value = -4728.4285
$ws.Range("H89").Value = 2447.9524
$ws.Range("I89").Value = 2430.7144
$ws.Range("J89").Value = 2482.4285
$ws.Range("K89").Value = 12153.572
$ws.Range("L89").Value = 12412.1425
$ws.Range("M89").Value = -6537.572
$ws.Range("N89").Value = -23644.1425
$ws.Range("H132").Value = 24407.773
$ws.Range("I132").Value = 1692.4
$ws.Range("J132").Value = 112745.336
$ws.Range("K132").Value = 5077.200000000001
$ws.Range("L132").Value = 338236.008
$ws.Range("M132").Value = -2547.200000000001
$ws.Range("N132").Value = -343296.008
$ws.Range("H134").Value = 37258.367
$ws.Range("I134").Value = 1145.1904
$ws.Range("J134").Value = 121522.445
$ws.Range("K134").Value = 3435.5712
$ws.Range("L134").Value = 364567.335
$ws.Range("M134").Value = -900.5711999999999
$ws.Range("N134").Value = -369637.335
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3214.7666
$ws.Range("I64").Value = 1493
$ws.Range("J64").Value = 3406.074
$ws.Range("K64").Value = 4479
$ws.Range("L64").Value = 10218.222
$ws.Range("M64").Value = -4209
$ws.Range("N64").Value = -10758.222
$ws.Range("H67").Value = 3214.7666
$ws.Range("I67").Value = 1493
$ws.Range("J67").Value = 3406.074
$ws.Range("K67").Value = 4479
$ws.Range("L67").Value = 10218.222
$ws.Range("M67").Value = -3543
$ws.Range("N67").Value = -12090.222
$ws.Range("H122").Value = 599.7273
$ws.Range("I122").Value = 282.6842
$ws.Range("K122").Value = 2544.1578
$ws.Range("M122").Value = -94.15779999999995
$ws.Range("H137").Value = 38937.438
$ws.Range("J137").Value = 47692.31
$ws.Range("L137").Value = 143076.93
$ws.Range("N137").Value = -153276.93
$ws.Range("H138").Value = 3429.6316
$ws.Range("I138").Value = 2203.75
$ws.Range("K138").Value = 6611.25
$ws.Range("M138").Value = -1471.25
$ws.Range("H141").Value = 8838.352999999999
$ws.Range("I141").Value = 5711.6924
$ws.Range("J141").Value = 19000
$ws.Range("K141").Value = 17135.0772
$ws.Range("L141").Value = 57000
$ws.Range("M141").Value = -11955.0772
$ws.Range("N141").Value = -67360
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 787.5
$ws.Range("I122").Value = 787.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2362.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 87.5
$ws.Range("N122").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2364.2334
$ws.Range("I40").Value = 2307.8276
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 2307.8276
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -2171.8276
$ws.Range("N40").Value = -4272
$ws.Range("H122").Value = 3221.853
$ws.Range("I122").Value = 2691.2942
$ws.Range("J122").Value = 3752.4119
$ws.Range("K122").Value = 8073.882599999999
$ws.Range("L122").Value = 11257.2357
$ws.Range("M122").Value = -5623.882599999999
$ws.Range("N122").Value = -16157.2357
$ws.Range("H132").Value = 54661.95
$ws.Range("I132").Value = 2160.75
$ws.Range("J132").Value = 334668.34
$ws.Range("K132").Value = 6482.25
$ws.Range("L132").Value = 1004005.02
$ws.Range("M132").Value = -3952.25
$ws.Range("N132").Value = -1009065.02

Write-Output "Applied all cell updates"